$wb = $excel.ActiveWorkbook

# ---- Sheet LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = 'Última actualización: 10:55:25'
$ws1.Cells.Item(3,1).Value = 'Total filas: 176'

$ws1.Cells.Item(37,1).Value = '06:33:46'
$ws1.Cells.Item(37,2).Value = '07:36'
$ws1.Cells.Item(37,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(37,4).Value = 63
$ws1.Cells.Item(37,5).Value = 'LP1912'

$ws1.Cells.Item(38,1).Value = '05:42:22'
$ws1.Cells.Item(38,2).Value = '07:36'
$ws1.Cells.Item(38,3).Value = '17X38_ROMERO'
$ws1.Cells.Item(38,4).Value = 114
$ws1.Cells.Item(38,5).Value = 'LP1912'

$ws1.Cells.Item(48,1).Value = '06:16:15'
$ws1.Cells.Item(48,2).Value = '08:00'
$ws1.Cells.Item(48,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(48,4).Value = 104
$ws1.Cells.Item(48,5).Value = 'LP1912'

$ws1.Cells.Item(49,1).Value = '06:33:46'
$ws1.Cells.Item(49,2).Value = '08:00'
$ws1.Cells.Item(49,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(49,4).Value = 87
$ws1.Cells.Item(49,5).Value = 'LP1912'

$ws1.Cells.Item(57,1).Value = '07:36:59'
$ws1.Cells.Item(57,2).Value = '08:14'
$ws1.Cells.Item(57,3).Value = '17_ROMERO'
$ws1.Cells.Item(57,4).Value = 38
$ws1.Cells.Item(57,5).Value = 'LP1912'

$ws1.Cells.Item(58,1).Value = '08:11:27'
$ws1.Cells.Item(58,2).Value = '08:14'
$ws1.Cells.Item(58,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(58,4).Value = 3
$ws1.Cells.Item(58,5).Value = 'LP1912'

$ws1.Cells.Item(59,1).Value = '07:48:35'
$ws1.Cells.Item(59,2).Value = '08:14'
$ws1.Cells.Item(59,3).Value = '10_OLMOS'
$ws1.Cells.Item(59,4).Value = 26
$ws1.Cells.Item(59,5).Value = 'LP1912'

$ws1.Cells.Item(140,1).Value = '10:55:25'
$ws1.Cells.Item(140,2).Value = '10:56'
$ws1.Cells.Item(140,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(140,4).Value = 1
$ws1.Cells.Item(140,5).Value = 'LP1912'

$ws1.Cells.Item(141,1).Value = '10:55:25'
$ws1.Cells.Item(141,2).Value = '10:56'
$ws1.Cells.Item(141,3).Value = '10_OLMOS'
$ws1.Cells.Item(141,4).Value = 1
$ws1.Cells.Item(141,5).Value = 'LP1912'

$ws1.Cells.Item(142,1).Value = '09:21:49'
$ws1.Cells.Item(142,2).Value = '10:56'
$ws1.Cells.Item(142,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(142,4).Value = 95
$ws1.Cells.Item(142,5).Value = 'LP1912'

$ws1.Cells.Item(143,1).Value = '10:04:17'
$ws1.Cells.Item(143,2).Value = '10:57'
$ws1.Cells.Item(143,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(143,4).Value = 53
$ws1.Cells.Item(143,5).Value = 'LP1912'

$ws1.Cells.Item(144,1).Value = '10:55:25'
$ws1.Cells.Item(144,2).Value = '10:59'
$ws1.Cells.Item(144,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(144,4).Value = 4
$ws1.Cells.Item(144,5).Value = 'LP1912'

$ws1.Cells.Item(145,1).Value = '09:21:49'
$ws1.Cells.Item(145,2).Value = '11:01'
$ws1.Cells.Item(145,3).Value = '17_ROMERO'
$ws1.Cells.Item(145,4).Value = 100
$ws1.Cells.Item(145,5).Value = 'LP1912'

$ws1.Cells.Item(146,1).Value = '10:36:18'
$ws1.Cells.Item(146,2).Value = '11:03'
$ws1.Cells.Item(146,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(146,4).Value = 27
$ws1.Cells.Item(146,5).Value = 'LP1912'

$ws1.Cells.Item(147,1).Value = '09:21:49'
$ws1.Cells.Item(147,2).Value = '11:04'
$ws1.Cells.Item(147,3).Value = '14_ABASTO'
$ws1.Cells.Item(147,4).Value = 103
$ws1.Cells.Item(147,5).Value = 'LP1912'

$ws1.Cells.Item(148,1).Value = '10:04:17'
$ws1.Cells.Item(148,2).Value = '11:05'
$ws1.Cells.Item(148,3).Value = '14_ABASTO'
$ws1.Cells.Item(148,4).Value = 61
$ws1.Cells.Item(148,5).Value = 'LP1912'

$ws1.Cells.Item(149,1).Value = '10:36:18'
$ws1.Cells.Item(149,2).Value = '11:11'
$ws1.Cells.Item(149,3).Value = '15_ABASTO'
$ws1.Cells.Item(149,4).Value = 35
$ws1.Cells.Item(149,5).Value = 'LP1912'

$ws1.Cells.Item(150,1).Value = '10:04:17'
$ws1.Cells.Item(150,2).Value = '11:11'
$ws1.Cells.Item(150,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(150,4).Value = 67
$ws1.Cells.Item(150,5).Value = 'LP1912'

$ws1.Cells.Item(151,1).Value = '09:21:49'
$ws1.Cells.Item(151,2).Value = '11:14'
$ws1.Cells.Item(151,3).Value = '225_C ROCA-H SUR'
$ws1.Cells.Item(151,4).Value = 113
$ws1.Cells.Item(151,5).Value = 'LP1912'

$ws1.Cells.Item(152,1).Value = '09:21:49'
$ws1.Cells.Item(152,2).Value = '11:20'
$ws1.Cells.Item(152,3).Value = '215C_EL PATO'
$ws1.Cells.Item(152,4).Value = 119
$ws1.Cells.Item(152,5).Value = 'LP1912'

$ws1.Cells.Item(153,1).Value = '10:04:17'
$ws1.Cells.Item(153,2).Value = '11:21'
$ws1.Cells.Item(153,3).Value = '215C_EL PATO'
$ws1.Cells.Item(153,4).Value = 77
$ws1.Cells.Item(153,5).Value = 'LP1912'

$ws1.Cells.Item(154,1).Value = '10:48:14'
$ws1.Cells.Item(154,2).Value = '11:21'
$ws1.Cells.Item(154,3).Value = '10_OLMOS'
$ws1.Cells.Item(154,4).Value = 33
$ws1.Cells.Item(154,5).Value = 'LP1912'

$ws1.Cells.Item(155,1).Value = '10:36:18'
$ws1.Cells.Item(155,2).Value = '11:22'
$ws1.Cells.Item(155,3).Value = '10_OLMOS'
$ws1.Cells.Item(155,4).Value = 46
$ws1.Cells.Item(155,5).Value = 'LP1912'

$ws1.Cells.Item(156,1).Value = '10:36:18'
$ws1.Cells.Item(156,2).Value = '11:24'
$ws1.Cells.Item(156,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(156,4).Value = 48
$ws1.Cells.Item(156,5).Value = 'LP1912'

$ws1.Cells.Item(157,1).Value = '10:36:18'
$ws1.Cells.Item(157,2).Value = '11:25'
$ws1.Cells.Item(157,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(157,4).Value = 49
$ws1.Cells.Item(157,5).Value = 'LP1912'

$ws1.Cells.Item(158,1).Value = '10:04:17'
$ws1.Cells.Item(158,2).Value = '11:25'
$ws1.Cells.Item(158,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(158,4).Value = 81
$ws1.Cells.Item(158,5).Value = 'LP1912'

$ws1.Cells.Item(159,1).Value = '10:04:17'
$ws1.Cells.Item(159,2).Value = '11:30'
$ws1.Cells.Item(159,3).Value = '15X38_ABASTO'
$ws1.Cells.Item(159,4).Value = 86
$ws1.Cells.Item(159,5).Value = 'LP1912'

$ws1.Cells.Item(160,1).Value = '10:48:14'
$ws1.Cells.Item(160,2).Value = '11:32'
$ws1.Cells.Item(160,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(160,4).Value = 44
$ws1.Cells.Item(160,5).Value = 'LP1912'

$ws1.Cells.Item(161,1).Value = '10:48:14'
$ws1.Cells.Item(161,2).Value = '11:33'
$ws1.Cells.Item(161,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(161,4).Value = 45
$ws1.Cells.Item(161,5).Value = 'LP1912'

$ws1.Cells.Item(162,1).Value = '10:36:18'
$ws1.Cells.Item(162,2).Value = '11:33'
$ws1.Cells.Item(162,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(162,4).Value = 57
$ws1.Cells.Item(162,5).Value = 'LP1912'

$ws1.Cells.Item(163,1).Value = '10:04:17'
$ws1.Cells.Item(163,2).Value = '11:34'
$ws1.Cells.Item(163,3).Value = '10_OLMOS'
$ws1.Cells.Item(163,4).Value = 90
$ws1.Cells.Item(163,5).Value = 'LP1912'

$ws1.Cells.Item(164,1).Value = '10:36:18'
$ws1.Cells.Item(164,2).Value = '11:35'
$ws1.Cells.Item(164,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(164,4).Value = 59
$ws1.Cells.Item(164,5).Value = 'LP1912'

$ws1.Cells.Item(165,1).Value = '10:04:17'
$ws1.Cells.Item(165,2).Value = '11:37'
$ws1.Cells.Item(165,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(165,4).Value = 93
$ws1.Cells.Item(165,5).Value = 'LP1912'

$ws1.Cells.Item(166,1).Value = '10:04:17'
$ws1.Cells.Item(166,2).Value = '11:40'
$ws1.Cells.Item(166,3).Value = '215A_EL PATO'
$ws1.Cells.Item(166,4).Value = 96
$ws1.Cells.Item(166,5).Value = 'LP1912'

$ws1.Cells.Item(167,1).Value = '10:55:25'
$ws1.Cells.Item(167,2).Value = '11:44'
$ws1.Cells.Item(167,3).Value = '215B_EL PATO'
$ws1.Cells.Item(167,4).Value = 49
$ws1.Cells.Item(167,5).Value = 'LP1912'

$ws1.Cells.Item(168,1).Value = '10:04:17'
$ws1.Cells.Item(168,2).Value = '11:45'
$ws1.Cells.Item(168,3).Value = '215B_EL PATO'
$ws1.Cells.Item(168,4).Value = 101
$ws1.Cells.Item(168,5).Value = 'LP1912'

$ws1.Cells.Item(169,1).Value = '10:55:25'
$ws1.Cells.Item(169,2).Value = '11:53'
$ws1.Cells.Item(169,3).Value = '15_ABASTO'
$ws1.Cells.Item(169,4).Value = 58
$ws1.Cells.Item(169,5).Value = 'LP1912'

$ws1.Cells.Item(170,1).Value = '10:04:17'
$ws1.Cells.Item(170,2).Value = '11:54'
$ws1.Cells.Item(170,3).Value = '225_GOMEZ'
$ws1.Cells.Item(170,4).Value = 110
$ws1.Cells.Item(170,5).Value = 'LP1912'

$ws1.Cells.Item(171,1).Value = '10:48:14'
$ws1.Cells.Item(171,2).Value = '12:07'
$ws1.Cells.Item(171,3).Value = '14_ABASTO'
$ws1.Cells.Item(171,4).Value = 79
$ws1.Cells.Item(171,5).Value = 'LP1912'

$ws1.Cells.Item(172,1).Value = '10:36:18'
$ws1.Cells.Item(172,2).Value = '12:29'
$ws1.Cells.Item(172,3).Value = '215C_EL PATO'
$ws1.Cells.Item(172,4).Value = 113
$ws1.Cells.Item(172,5).Value = 'LP1912'

$ws1.Cells.Item(173,1).Value = '10:36:18'
$ws1.Cells.Item(173,2).Value = '12:30'
$ws1.Cells.Item(173,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(173,4).Value = 114
$ws1.Cells.Item(173,5).Value = 'LP1912'

$ws1.Cells.Item(174,1).Value = '10:36:18'
$ws1.Cells.Item(174,2).Value = '12:31'
$ws1.Cells.Item(174,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(174,4).Value = 115
$ws1.Cells.Item(174,5).Value = 'LP1912'

$ws1.Cells.Item(175,1).Value = '10:48:14'
$ws1.Cells.Item(175,2).Value = '12:31'
$ws1.Cells.Item(175,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(175,4).Value = 103
$ws1.Cells.Item(175,5).Value = 'LP1912'

$ws1.Cells.Item(176,1).Value = '10:55:25'
$ws1.Cells.Item(176,2).Value = '12:36'
$ws1.Cells.Item(176,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(176,4).Value = 101
$ws1.Cells.Item(176,5).Value = 'LP1912'

$ws1.Cells.Item(177,1).Value = '10:48:14'
$ws1.Cells.Item(177,2).Value = '12:37'
$ws1.Cells.Item(177,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(177,4).Value = 109
$ws1.Cells.Item(177,5).Value = 'LP1912'

$ws1.Cells.Item(178,1).Value = '10:48:14'
$ws1.Cells.Item(178,2).Value = '12:40'
$ws1.Cells.Item(178,3).Value = '15X38_ABASTO'
$ws1.Cells.Item(178,4).Value = 112
$ws1.Cells.Item(178,5).Value = 'LP1912'

$ws1.Cells.Item(179,1).Value = '10:55:25'
$ws1.Cells.Item(179,2).Value = '12:42'
$ws1.Cells.Item(179,3).Value = '14_ABASTO'
$ws1.Cells.Item(179,4).Value = 107
$ws1.Cells.Item(179,5).Value = 'LP1912'

$ws1.Cells.Item(180,1).Value = '10:48:14'
$ws1.Cells.Item(180,2).Value = '12:43'
$ws1.Cells.Item(180,3).Value = '14_ABASTO'
$ws1.Cells.Item(180,4).Value = 115
$ws1.Cells.Item(180,5).Value = 'LP1912'

$ws1.Cells.Item(181,1).Value = '10:55:25'
$ws1.Cells.Item(181,2).Value = '12:43'
$ws1.Cells.Item(181,3).Value = '15X38_ABASTO'
$ws1.Cells.Item(181,4).Value = 108
$ws1.Cells.Item(181,5).Value = 'LP1912'

# ---- Sheet LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = 'Última actualización: 10:55:25'
$ws2.Cells.Item(3,1).Value = 'Total filas: 25'

$ws2.Cells.Item(28,1).Value = '10:55:25'
$ws2.Cells.Item(28,2).Value = '11:44'
$ws2.Cells.Item(28,3).Value = '215B_EL PATO'
$ws2.Cells.Item(28,4).Value = 49
$ws2.Cells.Item(28,5).Value = 'LP1912'

$ws2.Cells.Item(29,1).Value = '10:04:17'
$ws2.Cells.Item(29,2).Value = '11:45'
$ws2.Cells.Item(29,3).Value = '215B_EL PATO'
$ws2.Cells.Item(29,4).Value = 101
$ws2.Cells.Item(29,5).Value = 'LP1912'

$ws2.Cells.Item(30,1).Value = '10:36:18'
$ws2.Cells.Item(30,2).Value = '12:29'
$ws2.Cells.Item(30,3).Value = '215C_EL PATO'
$ws2.Cells.Item(30,4).Value = 113
$ws2.Cells.Item(30,5).Value = 'LP1912'

# ---- Sheet 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = 'Última actualización: 10:55:25'
